$wb = $excel.ActiveWorkbook

# --- Update Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 1.8.1 -> 1.8.2
$meta.Range("B3").Value = "1.8.2"

# Date: 2023-06-27T22:42:19-04:00 -> 2023-09-01T14:45:29-04:00
$meta.Range("B8").Value = "2023-09-01T14:45:29-04:00"

# --- Update Elements sheet ---
$elements = $wb.Worksheets.Item("Elements")

# Add the ele-1 / ext-1 invariant text to the top-level Extension row (row 1, column AJ)
$elements.Range("AJ1").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}" + [char]10 + "ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

$wb.Save()
